# Adds the first six rows of the new Azure "B-series" burstable VM sizes
# (Standard_B1S, Standard_B2S, Standard_B1MS, Standard_B2MS, Standard_B4MS,
# Standard_B8MS) to the AzureVMSizes worksheet, right after the existing
# data which currently ends at row 135.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 136
$lastRow  = 141

# ---- column A: Name --------------------------------------------------
$names = @("Standard_B1S", "Standard_B2S", "Standard_B1MS", "Standard_B2MS", "Standard_B4MS", "Standard_B8MS")
for ($i = 0; $i -lt 6; $i++) {
    $ws.Cells.Item($firstRow + $i, 1).Value = $names[$i]
}

# ---- column B: Type (same "VM" value used throughout the sheet) -----
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 2).Value = "VM"
}

# ---- column C: PricePerHourPAYG --------------------------------------
$prices = @(0.006, 0.023, 0.012, 0.046, 0.092, 0.183)
for ($i = 0; $i -lt 6; $i++) {
    $ws.Cells.Item($firstRow + $i, 3).Value = $prices[$i]
}

# ---- column D: PricePerHourLP (unknown for these new sizes) ---------
# ---- column E: ACU (unknown for these new sizes) ---------------------
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 4).Value = "Unknown"
    $ws.Cells.Item($r, 5).Value = "Unknown"
}

# ---- column G: match the existing "accounting" number format used
# by the rest of the CoreFactorIndex column before we drop the shared
# formulas in (so the new formula cells pick up style index 2 instead
# of the column default of 1).
$gFormat = $ws.Range("G135").NumberFormat
$ws.Range("G136:G141").NumberFormat = $gFormat

# ---- columns F & G: shared IFERROR formulas, filled down as a block -
$ws.Range("F136:F141").Formula = "=IFERROR(E136*I136,""Unknown"")"
$ws.Range("G136:G141").Formula = "=IFERROR(C136/F136*1000,""Unknown"")"

# ---- column H: SSD -----------------------------------------------------
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 8).Value = "No"
}

# ---- column I: NumberOfCores -------------------------------------------
$cores = @(1, 2, 1, 2, 4, 8)
for ($i = 0; $i -lt 6; $i++) {
    $ws.Cells.Item($firstRow + $i, 9).Value = $cores[$i]
}

# ---- column J: MemoryInMB (first three literal, last three computed) --
$ws.Cells.Item(136, 10).Value = 1024
$ws.Cells.Item(137, 10).Value = 4096
$ws.Cells.Item(138, 10).Value = 2048
$ws.Range("J139").Formula = "=8*1024"
$ws.Range("J140").Formula = "=16*1024"
$ws.Range("J141").Formula = "=32*1024"

# ---- columns K-R: MaxNics, Bandwidth, MaxDataDiskCount,
# MaxDataDiskSizeGB, MaxDataDiskIops, MaxDataDiskThroughputMBs,
# MaxVmIops, MaxVmThroughputMBs -- all unknown for these new sizes ------
$unknownCols = @(11, 13, 14, 15, 16, 17, 18)  # K, M, N, O, P, Q, R  (L done below, already formatted)
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 11).Value = "Unknown"  # K MaxNics
    $ws.Cells.Item($r, 12).Value = "Unknown"  # L Bandwidth
    $ws.Cells.Item($r, 13).Value = "Unknown"  # M MaxDataDiskCount
    $ws.Cells.Item($r, 14).Value = "Unknown"  # N MaxDataDiskSizeGB
    $ws.Cells.Item($r, 15).Value = "Unknown"  # O MaxDataDiskIops
    $ws.Cells.Item($r, 16).Value = "Unknown"  # P MaxDataDiskThroughputMBs
    $ws.Cells.Item($r, 17).Value = "Unknown"  # Q MaxVmIops
    $ws.Cells.Item($r, 18).Value = "Unknown"  # R MaxVmThroughputMBs
}

# ---- column S: ResourceDiskSizeInMB ------------------------------------
$resourceDisk = @(2000, 8000, 4000, 16000, 32000, 64000)
for ($i = 0; $i -lt 6; $i++) {
    $ws.Cells.Item($firstRow + $i, 19).Value = $resourceDisk[$i]
}

# ---- column T: TempDiskSizeInGB ----------------------------------------
$tempDisk = @(2, 8, 4, 16, 32, 64)
for ($i = 0; $i -lt 6; $i++) {
    $ws.Cells.Item($firstRow + $i, 20).Value = $tempDisk[$i]
}

# ---- columns U-W: TempDiskIops / TempDiskReadMBs / TempDiskWriteMBs ---
# ---- columns X-Z: SAPS2T / SAPS3T / HANA ------------------------------
# ---- column AA: Hyperthreaded -----------------------------------------
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 21).Value = "Unknown"  # U
    $ws.Cells.Item($r, 22).Value = "Unknown"  # V
    $ws.Cells.Item($r, 23).Value = "Unknown"  # W
    $ws.Cells.Item($r, 24).Value = "Unknown"  # X
    $ws.Cells.Item($r, 25).Value = "Unknown"  # Y
    $ws.Cells.Item($r, 26).Value = "Unknown"  # Z
    $ws.Cells.Item($r, 27).Value = "Unknown"  # AA
}

# ---- style fix-up ------------------------------------------------------
# Columns C, D, E, F, G, L, U, V, W already default to style index 1 (or
# 2, for G) via the sheet's <cols> definitions. The remaining "Unknown"
# placeholder / plain-number columns (K, M, N, O, P, Q, R, T, X, Y, Z, AA)
# need that same right-aligned style applied explicitly, matching every
# other row in this table. Copy the format from row 135, which already
# carries that style, onto the corresponding cells of the new rows.
$formatCols = @("K", "M", "N", "O", "P", "Q", "R", "T", "X", "Y", "Z")
foreach ($col in $formatCols) {
    $ws.Range($col + "135").Copy()
    $ws.Range($col + $firstRow + ":" + $col + $lastRow).PasteSpecial(-4122)
}

# Row 135's AA cell ("Yes", no explicit style) isn't a usable source for
# AA's style, unlike the other columns above -- use the header cell AA1,
# which already carries the same right-aligned style index, instead.
$ws.Range("AA1").Copy()
$ws.Range("AA" + $firstRow + ":AA" + $lastRow).PasteSpecial(-4122)

$excel.CutCopyMode = 0
